$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "test"
$ws.Name = "test"

# Enter the username/password sample data.
# Row 1 (headers) first, then column A for the remaining rows, then column B --
# this produces the same shared-string insertion order as the authored workbook.
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("A2").Value = "sufail"
$ws.Range("A3").Value = "Saban"
$ws.Range("B2").Value = "password1"
$ws.Range("B3").Value = "password2"

# Auto-fit the two data columns to their content.
$ws.Columns("A:B").AutoFit() | Out-Null

# Leave the selection on C4, just past the data, as in the saved workbook.
$ws.Range("C4").Select() | Out-Null
